$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" column (column G) values recomputed from the regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals"). Only column G (rows 2-21) changes.
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 3
    6  = 4
    7  = 1
    8  = 3
    9  = 10
    10 = 9
    11 = 6
    12 = 4
    13 = 6
    14 = 2
    15 = 10
    16 = 7
    17 = 3
    18 = 5
    19 = 3
    20 = 5
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
